# Update the users and fixture list:
# add four new Euro 2024 fixtures (Spain vs Germany, Portugal vs France,
# England vs Switzerland, Netherlands vs Turkey) as new header columns
# AV:AY, and fill in the already-known results for the users who picked
# them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) needs the same bold/bordered/centered style as the
# rest of the header cells (e.g. AU1) - copy formats from AU1 first, then
# set the new header captions.
$ws.Range("AU1").Copy()
$ws.Range("AV1:AY1").PasteSpecial(-4122)

$ws.Range("AV1").Value = "Spain vs Germany"
$ws.Range("AW1").Value = "Portugal vs France"
$ws.Range("AX1").Value = "England vs Switzerland"
$ws.Range("AY1").Value = "Netherlands vs Turkey"

# Row 2 results
$ws.Range("AV2").Value = "['Spain', 3, 2]"
$ws.Range("AW2").Value = "['France', 2, 3]"
$ws.Range("AX2").Value = "['Switzerland', 1, 2]"
$ws.Range("AY2").Value = "['Netherlands', 3, 2]"

# Row 5 results
$ws.Range("AV5").Value = "['Draw', 1, 1]"
$ws.Range("AW5").Value = "['Draw', 1, 1]"
$ws.Range("AX5").Value = "['England', 1, 0]"
$ws.Range("AY5").Value = "['Netherlands', 2, 1]"

# Row 7 results
$ws.Range("AV7").Value = "['Germany', 1, 2]"
$ws.Range("AW7").Value = "['Draw', 1, 1]"
$ws.Range("AX7").Value = "['Switzerland', 1, 2]"
$ws.Range("AY7").Value = "['Turkey', 1, 2]"

# Row 13 results
$ws.Range("AV13").Value = "['Spain', 2, 1]"
$ws.Range("AW13").Value = "['Portugal', 2, 1]"
$ws.Range("AX13").Value = "['England', 1, 0]"
$ws.Range("AY13").Value = "['Draw', 1, 1]"
